# Add version numbers to schemas.
#
# 1) Insert a new "version" column at the front of the "Export as TSV"
#    header row (shifting affiliation/first_name/last_name/
#    middle_name_or_initial/name/orcid_id one column to the right, along
#    with their header comments).
# 2) Add a data validation on column A (rows 2+) restricting values to a
#    lookup list.
# 3) Add a new "version list" worksheet holding the allowed version
#    number(s) that the validation list formula points at.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) Shift the header row's values one column to the right (F->G,
#        E->F, ... A->B) then set the new "version" header in A1.
$ws1.Range("G1").Value2 = $ws1.Range("F1").Value2
$ws1.Range("F1").Value2 = $ws1.Range("E1").Value2
$ws1.Range("E1").Value2 = $ws1.Range("D1").Value2
$ws1.Range("D1").Value2 = $ws1.Range("C1").Value2
$ws1.Range("C1").Value2 = $ws1.Range("B1").Value2
$ws1.Range("B1").Value2 = $ws1.Range("A1").Value2
$ws1.Range("A1").Value2 = "version"

# G1 is a brand-new cell - give it the same bold/centered/wrap header
# styling the other header cells already carry (style copied from F1).
$ws1.Range("G1").Font.Bold = $true
$ws1.Range("G1").HorizontalAlignment = -4108  # xlCenter
$ws1.Range("G1").WrapText = $true

# --- Shift the header cell-comments the same way, then set A1's new
#     comment text.
$oldF = $ws1.Range("F1").Comment.Text()
$oldE = $ws1.Range("E1").Comment.Text()
$oldD = $ws1.Range("D1").Comment.Text()
$oldC = $ws1.Range("C1").Comment.Text()
$oldB = $ws1.Range("B1").Comment.Text()
$oldA = $ws1.Range("A1").Comment.Text()

$ws1.Range("G1").AddComment($oldF) | Out-Null
$ws1.Range("F1").Comment.Text($oldE) | Out-Null
$ws1.Range("E1").Comment.Text($oldD) | Out-Null
$ws1.Range("D1").Comment.Text($oldC) | Out-Null
$ws1.Range("C1").Comment.Text($oldB) | Out-Null
$ws1.Range("B1").Comment.Text($oldA) | Out-Null
$ws1.Range("A1").Comment.Text("Version of the schema to use when validating this metadata.") | Out-Null

# --- 2) Add the new "version list" worksheet after the existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "version list"

# Populate A1 with the text value "1" (stored as a real string, not a
# number) - build it via a formula then freeze the result to a value so
# no formula/style residue is left behind.
$ws2.Range("A1").Formula = "=""1"""
$ws2.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# --- 3) Restrict column A (below the header) to values from the
#     "version list" sheet.
$rng = $ws1.Range("A2:A1048576")
$rng.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")  # xlValidateList
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: 1."
$rng.Validation.ShowInput = $true
$rng.Validation.ShowError = $true

# Leave the original sheet as the active/selected tab, matching the
# workbook's prior view state.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
